# Apply data updates to the "Resumo Inscrições Integrado" sheet
# per commit "Data update using git"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 103
$ws.Range("E4").Value = 48
$ws.Range("E5").Value = 144
$ws.Range("E9").Value = 9
$ws.Range("E10").Value = 614
$ws.Range("F10").Value = 316
$ws.Range("H10").Value = 412
$ws.Range("E11").Value = 397
$ws.Range("F11").Value = 221
$ws.Range("H11").Value = 285
$ws.Range("E12").Value = 608
$ws.Range("F12").Value = 343
$ws.Range("H12").Value = 429
$ws.Range("E13").Value = 146
$ws.Range("F13").Value = 80
$ws.Range("H13").Value = 114
$ws.Range("E16").Value = 213
$ws.Range("E17").Value = 110
$ws.Range("F17").Value = 61
$ws.Range("H17").Value = 85
$ws.Range("E18").Value = 55
$ws.Range("E20").Value = 94
$ws.Range("E21").Value = 144
$ws.Range("E22").Value = 179
$ws.Range("E23").Value = 211
$ws.Range("E24").Value = 233
$ws.Range("F24").Value = 128
$ws.Range("H24").Value = 158
$ws.Range("E25").Value = 294
$ws.Range("F25").Value = 156
$ws.Range("H25").Value = 216
$ws.Range("E26").Value = 166
$ws.Range("E27").Value = 349
$ws.Range("E28").Value = 210
$ws.Range("E29").Value = 176
$ws.Range("F29").Value = 104
$ws.Range("H29").Value = 145
$ws.Range("E30").Value = 228
$ws.Range("E31").Value = 76
$ws.Range("E33").Value = 310
$ws.Range("F33").Value = 163
$ws.Range("G33").Value = 90
$ws.Range("E34").Value = 233
$ws.Range("F34").Value = 159
$ws.Range("G34").Value = 38
$ws.Range("H34").Value = 197
$ws.Range("E36").Value = 80
$ws.Range("F36").Value = 50
$ws.Range("H36").Value = 60
$ws.Range("E37").Value = 176
$ws.Range("F37").Value = 91
$ws.Range("H37").Value = 128
$ws.Range("E38").Value = 97
$ws.Range("E39").Value = 188
$ws.Range("E40").Value = 280
$ws.Range("F40").Value = 132
$ws.Range("H40").Value = 212
$ws.Range("E42").Value = 406
$ws.Range("F42").Value = 231
$ws.Range("H42").Value = 292
$ws.Range("E43").Value = 131
$ws.Range("F43").Value = 70
$ws.Range("H43").Value = 97
$ws.Range("E44").Value = 326
$ws.Range("E46").Value = 351
$ws.Range("F46").Value = 193
$ws.Range("H46").Value = 256
$ws.Range("E47").Value = 490
$ws.Range("F47").Value = 259
$ws.Range("H47").Value = 351
$ws.Range("E48").Value = 233
$ws.Range("E49").Value = 305
$ws.Range("F49").Value = 141
$ws.Range("H49").Value = 228
$ws.Range("E50").Value = 253
$ws.Range("E51").Value = 249
$ws.Range("F51").Value = 117
$ws.Range("H51").Value = 191
$ws.Range("E52").Value = 31

$wb.Save()
